$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Label the AVERAGE/СРЗНАЧ formula used in B13 (adds a new shared string)
$ws.Range("F13").Value = "ФОРМУЛА =СРЗНАЧ(B2:B9)"

# Label the discounted-revenue formula used in column E (adds a new shared string)
$ws.Range("H2").Value = "ФОРМУЛА =(B2-(B2*C2))*D2"

# Leave the selection where the user ended up after typing the labels
$ws.Range("H3").Select()
